$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 7
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "2023-11-22 | Jan:01 - препарат принят (От поноса)"
$ws.Range("E9").Value = "22/11/2023 20:09:39"
